$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column C: Reason header + two reason values
$ws.Range("C1").Value = "Reason"
$ws.Range("C2").Value = "Inappropriate content"
$ws.Range("C3").Value = "Fraud Email"

# Update the selected cell to match the target (E13)
$ws.Range("E13").Select()
